# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Columns E (nivel-formativo-grupo-iaest-descripcion), F (aragon) and H (sexo)
# are re-classified: E and H move from "dimension" to "measure", and F moves
# from a curated skos:Concept dimension with its own mapping file to a plain
# sdmx-dimension:refArea / URI-Comunidad dimension (like the other refArea
# columns C and J). The old per-column mapping-file row (row 5) is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimension/measure qualifier
$ws.Range("E2").Value = "iaest-measure:nivel-formativo-grupo-iaest-descripcion"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "iaest-measure:sexo"

# Row 3: medida/dim label
$ws.Range("E3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Row 4: datatype / URI kind
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("H4").Value = "xsd:int"

# Row 5 (per-column mapping workbook references) no longer applies - remove it
$ws.Rows(5).Delete()
